# Update TPM-derived NATMI metrics (Cd38-Pecam1) to new values per updated scripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.36822333333333
$ws.Range("H2").Value = 136.10467
$ws.Range("I2").Value = 0.9399475154355185
$ws.Range("J2").Value = 0.9399475154355185
$ws.Range("M2").Value = 186.2750726666667
$ws.Range("N2").Value = 558.8252179999999
$ws.Range("O2").Value = 0.9729431886728379
$ws.Range("P2").Value = 0.9729431886728379
$ws.Range("Q2").Value = 8450.969098174228
$ws.Range("R2").Value = 76058.72188356805
$ws.Range("S2").Value = 0.9145155328529448
$ws.Range("T2").Value = 0.9145155328529448

$ws.Range("G3").Value = 45.36822333333333
$ws.Range("H3").Value = 136.10467
$ws.Range("I3").Value = 0.9399475154355185
$ws.Range("J3").Value = 0.9399475154355185
$ws.Range("O3").Value = 0.01102259370028598
$ws.Range("P3").Value = 0.01102259370028598
$ws.Range("Q3").Value = 95.74207397444444
$ws.Range("R3").Value = 861.67866577
$ws.Range("S3").Value = 0.010360659562239
$ws.Range("T3").Value = 0.010360659562239

$ws.Range("G4").Value = 45.36822333333333
$ws.Range("H4").Value = 136.10467
$ws.Range("I4").Value = 0.9399475154355185
$ws.Range("J4").Value = 0.9399475154355185
$ws.Range("M4").Value = 3.069835
$ws.Range("N4").Value = 9.209505
$ws.Range("O4").Value = 0.01603421762687604
$ws.Range("P4").Value = 0.01603421762687604
$ws.Range("Q4").Value = 139.2729598764833
$ws.Range("R4").Value = 1253.45663888835
$ws.Range("S4").Value = 0.01507132302033453
$ws.Range("T4").Value = 0.01507132302033453

$ws.Range("I5").Value = 0.05085382185741109
$ws.Range("J5").Value = 0.05085382185741109
$ws.Range("M5").Value = 186.2750726666667
$ws.Range("N5").Value = 558.8252179999999
$ws.Range("O5").Value = 0.9729431886728379
$ws.Range("P5").Value = 0.9729431886728379
$ws.Range("Q5").Value = 457.2213554305848
$ws.Range("R5").Value = 4114.992198875263
$ws.Range("S5").Value = 0.04947787959415001
$ws.Range("T5").Value = 0.04947787959415001

$ws.Range("I6").Value = 0.05085382185741109
$ws.Range("J6").Value = 0.05085382185741109
$ws.Range("O6").Value = 0.01102259370028598
$ws.Range("P6").Value = 0.01102259370028598
$ws.Range("Q6").Value = 5.179917276444443
$ws.Range("R6").Value = 46.61925548799999
$ws.Range("S6").Value = 0.0005605410164409648
$ws.Range("T6").Value = 0.0005605410164409648

$ws.Range("I7").Value = 0.05085382185741109
$ws.Range("J7").Value = 0.05085382185741109
$ws.Range("M7").Value = 3.069835
$ws.Range("N7").Value = 9.209505
$ws.Range("O7").Value = 0.01603421762687604
$ws.Range("P7").Value = 0.01603421762687604
$ws.Range("Q7").Value = 7.535061452693332
$ws.Range("R7").Value = 67.81555307424
$ws.Range("S7").Value = 0.0008154012468201149
$ws.Range("T7").Value = 0.0008154012468201149

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.4439896666666667
$ws.Range("H8").Value = 1.331969
$ws.Range("I8").Value = 0.009198662707070463
$ws.Range("J8").Value = 0.009198662707070463
$ws.Range("M8").Value = 186.2750726666667
$ws.Range("N8").Value = 558.8252179999999
$ws.Range("O8").Value = 0.9729431886728379
$ws.Range("P8").Value = 0.9729431886728379
$ws.Range("Q8").Value = 82.70420742158244
$ws.Range("R8").Value = 744.3378667942419
$ws.Range("S8").Value = 0.008949776225743056
$ws.Range("T8").Value = 0.008949776225743056

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.4439896666666667
$ws.Range("H9").Value = 1.331969
$ws.Range("I9").Value = 0.009198662707070463
$ws.Range("J9").Value = 0.009198662707070463
$ws.Range("O9").Value = 0.01102259370028598
$ws.Range("P9").Value = 0.01102259370028598
$ws.Range("Q9").Value = 0.9369661932222222
$ws.Range("R9").Value = 8.432695739
$ws.Range("S9").Value = 0.0001013931216060104
$ws.Range("T9").Value = 0.0001013931216060104

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4439896666666667
$ws.Range("H10").Value = 1.331969
$ws.Range("I10").Value = 0.009198662707070463
$ws.Range("J10").Value = 0.009198662707070463
$ws.Range("M10").Value = 3.069835
$ws.Range("N10").Value = 9.209505
$ws.Range("O10").Value = 0.01603421762687604
$ws.Range("P10").Value = 0.01603421762687604
$ws.Range("Q10").Value = 1.362975018371667
$ws.Range("R10").Value = 12.266775165345
$ws.Range("S10").Value = 0.0001474933597213964
$ws.Range("T10").Value = 0.0001474933597213964

